$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that used to sit right after
#    "wykonuje backupy?".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Drop the "[Kod QR? kod kreskowy?]" note after "wbijam produkty", leaving
#    just the trailing space that was already in front of it.
$d.Content.Find.Execute(" [Kod QR? kod kreskowy?]", $true, $false, $false, `
    $false, $false, $true, 1, $false, " ", 2) | Out-Null

# 3. Remove the whole "Android do skanowania kodów" bullet point (the entire
#    paragraph, including its paragraph mark/numbering).
$rng = $d.Content
$found = $rng.Find.Execute("Android do skanowania kodów", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Paragraphs(1).Range.Delete()
}

# 4. Word's "_GoBack" bookmark now belongs at the start of the next edit
#    point, which is right before "Diagram klas:".
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Diagram klas:", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $target = $rng2.Paragraphs(1).Range.Duplicate
    $target.Collapse(1)
    $d.Bookmarks.Add("_GoBack", $target)
}
